# Apply the commit "Refined metadata to be additional tab":
#  1. Update the "panel_query_time" timestamps in column F of the "data" sheet
#     (the panel data was re-queried at a later time).
#  2. Add a new "metadata" worksheet (after "data") describing the panel query
#     that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Refresh the time_taken / panel_query_time values on the "data" sheet
# ---------------------------------------------------------------------------
$newTimes = @(
  "2021-10-05 14:35:43.390154",
  "2021-10-05 14:35:43.390161",
  "2021-10-05 14:35:43.390165",
  "2021-10-05 14:35:43.390168",
  "2021-10-05 14:35:43.390171",
  "2021-10-05 14:35:43.390173",
  "2021-10-05 14:35:43.390176",
  "2021-10-05 14:35:43.390179",
  "2021-10-05 14:35:43.390182",
  "2021-10-05 14:35:43.390184",
  "2021-10-05 14:35:43.390187",
  "2021-10-05 14:35:43.390189",
  "2021-10-05 14:35:43.390192",
  "2021-10-05 14:35:43.390194",
  "2021-10-05 14:35:43.390197",
  "2021-10-05 14:35:43.390200",
  "2021-10-05 14:35:43.390202",
  "2021-10-05 14:35:43.390205",
  "2021-10-05 14:35:43.390208",
  "2021-10-05 14:35:43.390210",
  "2021-10-05 14:35:43.390213",
  "2021-10-05 14:35:43.390215",
  "2021-10-05 14:35:43.390218",
  "2021-10-05 14:35:43.390220",
  "2021-10-05 14:35:43.390223",
  "2021-10-05 14:35:43.390226",
  "2021-10-05 14:35:43.390229",
  "2021-10-05 14:35:43.390231",
  "2021-10-05 14:35:43.390234",
  "2021-10-05 14:35:43.390236",
  "2021-10-05 14:35:43.390239",
  "2021-10-05 14:35:43.390241",
  "2021-10-05 14:35:43.390245",
  "2021-10-05 14:35:43.390247"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
  $dataSheet.Cells.Item($i + 2, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the "metadata" worksheet right after "data"
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Copy the header / index cell formatting from the "data" sheet so the new
# sheet's header row (B1:G1) and index cell (A2) use the same bold/bordered
# style (cell style index 1) instead of creating brand-new styles. Copy one
# column at a time so there is no ambiguity from range-size "tiling".
foreach ($col in 2..7) {
  $dataSheet.Cells.Item(1, 2).Copy()
  $meta.Cells.Item(1, $col).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$dataSheet.Cells.Item(2, 1).Copy()
$meta.Cells.Item(2, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wb.Application.CutCopyMode = 0

# Header row
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row (pandas-style integer index 0 in column A)
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Short Rib Polydactyly_Jeune Asphyxiating Thoracic Dystrophy_Skeletal Ciliopathy"
$meta.Cells.Item(2, 3).Value = 179

# Force "1.5" to be stored as text (matches source data, not a numeric 1.5)
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.5"

$meta.Cells.Item(2, 5).Value = "2021-08-12T23:34:19.399737Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:43.386392"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/179/?format=json"

# Leave the "data" sheet active/selected as it was originally.
$dataSheet.Activate()
